$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44326
$ws.Range("M2").Value = 65
$ws.Range("D3").Value = 44326
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("S3").Value = 800
$ws.Range("D4").Value = 44314
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 47
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 9000
$ws.Range("P4").Value = 9000
$ws.Range("S4").Value = 900
$ws.Range("D5").Value = 44329
$ws.Range("M5").Value = 56
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("D6").Value = 44329
$ws.Range("M6").Value = 50
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("D7").Value = 44307
$ws.Range("M7").Value = 40
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 44319
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 57
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 800
$ws.Range("D10").Value = 44322
$ws.Range("M10").Value = 56
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 800
$ws.Range("D12").Value = 44302
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 1000
$ws.Range("D13").Value = 44343
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 47
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("D14").Value = 44343
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 9000
$ws.Range("P14").Value = 9000
$ws.Range("R14").Value = "Región Metropolitana"
$ws.Range("S14").Value = 900
$ws.Range("D15").Value = 44343
$ws.Range("M15").Value = 58
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 800
$ws.Range("D17").Value = 44328
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 45
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 800
$ws.Range("D18").Value = 44328
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 48
$ws.Range("N18").Value = 7000
$ws.Range("O18").Value = 7000
$ws.Range("P18").Value = 7000
$ws.Range("S18").Value = 700
$ws.Range("D19").Value = 44308
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 48
$ws.Range("N20").Value = 8000
$ws.Range("O20").Value = 8000
$ws.Range("P20").Value = 8000
$ws.Range("S20").Value = 800
$ws.Range("D21").Value = 44315
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("S21").Value = 1000
$ws.Range("D22").Value = 44321
$ws.Range("M22").Value = 58
$ws.Range("N22").Value = 9000
$ws.Range("O22").Value = 9000
$ws.Range("P22").Value = 9000
$ws.Range("S22").Value = 900
$ws.Range("D23").Value = 44306
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("S23").Value = 1000
$ws.Range("D24").Value = 44699
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 56
$ws.Range("N24").Value = 12000
$ws.Range("O24").Value = 12000
$ws.Range("P24").Value = 12000
$ws.Range("S24").Value = 1200
$ws.Range("D25").Value = 44699
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 60
$ws.Range("R25").Value = "Provincia de Quillota"
$ws.Range("D26").Value = 44333
$ws.Range("L26").Value = "Especial"
$ws.Range("M26").Value = 58
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("R26").Value = "Provincia de Quillota"
$ws.Range("S26").Value = 1000
$ws.Range("D27").Value = 44333
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 65
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("R27").Value = "Provincia de Quillota"
$ws.Range("S27").Value = 900
$ws.Range("D28").Value = 44333
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("S28").Value = 800
$ws.Range("D29").Value = 44301
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 45
$ws.Range("N29").Value = 10000
$ws.Range("O29").Value = 10000
$ws.Range("P29").Value = 10000
$ws.Range("S29").Value = 1000
$ws.Range("D30").Value = 44323
$ws.Range("M30").Value = 60
$ws.Range("D31").Value = 44323
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 50
$ws.Range("R31").Value = "Provincia de Quillota"
$ws.Range("D32").Value = 44309
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 45
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 10000
$ws.Range("R32").Value = "Provincia de Quillota"
$ws.Range("S32").Value = 1000
